# Applies the "processing" -> "estimator" terminology change on the
# mitigation-stages diagram and updates the fixed date footer text from
# 12/29/2020 to 2/8/2021 everywhere it appears (the Slide Master and every
# Custom Layout that carries a Date placeholder).

$p = $ppt.ActivePresentation

$oldDate = "12/29/2020"
$newDate = "2/8/2021"
$ppPlaceholderDate = 16

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $ph = $sh.PlaceholderFormat
            if ($ph.Type -eq $ppPlaceholderDate) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide Master date placeholder.
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every Custom Layout's date placeholder.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DateShapes $layout.Shapes
}

# Slide content: rename the three "processing" stage labels to the new
# "estimator mitigation" terminology.
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "preprocessing") {
            $tr.Text = "pre-estimator mitigation"
        }
        elseif ($tr.Text -eq "inprocessing") {
            $tr.Text = "in-estimator mitigation"
        }
        elseif ($tr.Text -eq "postprocessing") {
            $tr.Text = "post-estimator mitigation"
        }
    }
}
